$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows being appended to the report (update through 6/03).
$newRows = @(
    @(245, 44319, 0, 17, 519.0839694656488),
    @(246, 44320, 0, 17, 519.0839694656488),
    @(247, 44321, 0, 16, 488.5496183206107)
)

$lastRow = 244

foreach ($r in $newRows) {
    $rowNum = $r[0]

    # Copy the formatting (date style with border/alignment) from the last
    # existing row in column A so the new date cell matches the rest of
    # the column.
    $ws.Range("A$lastRow").Copy() | Out-Null
    $ws.Range("A$rowNum").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}

$excel.CutCopyMode = 0
